$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G4").Value = 99
$ws.Range("H4").Value = 1145
$ws.Range("I4").Value = 1012
$ws.Range("J4").Value = 1101
$ws.Range("Q4").Value = 734
